$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# Row 6: CC_cal_temp - numeric value
$ws.Range("H6").Value = 20.06

# Row 7: CC_wl - array of wavelengths (stored as text)
$ws.Range("H7").Value = "[190.45, 191.24, 192.03, 192.82, 193.61, 194.4, 195.19, 195.99, 196.78, 197.57, 198.36, 199.15, 199.95, 200.74, 201.53, 202.32, 203.12, 203.91, 204.7, 205.5, 206.29, 207.09, 207.88, 208.67, 209.47, 210.26, 211.06, 211.85, 212.65, 213.44, 214.24, 215.04, 215.83, 216.63, 217.42, 218.22, 219.02, 219.81, 220.61, 221.41, 222.2, 223, 223.8, 224.6, 225.39, 226.19, 226.99, 227.79, 228.59, 229.38, 230.18, 230.98, 231.78, 232.58, 233.38, 234.18, 234.97, 235.77, 236.57, 237.37, 238.17, 238.97, 239.77, 240.57, 241.37, 242.17, 242.97, 243.77, 244.57, 245.37, 246.17, 246.97, 247.77, 248.57, 249.37, 250.17, 250.98, 251.78, 252.58, 253.38, 254.18, 254.98, 255.78, 256.58, 257.39, 258.19, 258.99, 259.79, 260.59, 261.39, 262.19, 263, 263.8, 264.6, 265.4, 266.2, 267.01, 267.81, 268.61, 269.41, 270.21, 271.02, 271.82, 272.62, 273.42, 274.23, 275.03, 275.83, 276.63, 277.44, 278.24, 279.04, 279.84, 280.65, 281.45, 282.25, 283.05, 283.85, 284.66, 285.46, 286.26, 287.06, 287.87, 288.67, 289.47, 290.27, 291.08, 291.88, 292.68, 293.48, 294.29, 295.09, 295.89, 296.69, 297.49, 298.3, 299.1, 299.9, 300.7, 301.51, 302.31, 303.11, 303.91, 304.71, 305.51, 306.32, 307.12, 307.92, 308.72, 309.52, 310.32, 311.13, 311.93, 312.73, 313.53, 314.33, 315.13, 315.93, 316.73, 317.53, 318.34, 319.14, 319.94, 320.74, 321.54, 322.34, 323.14, 323.94, 324.74, 325.54, 326.34, 327.14, 327.94, 328.74, 329.54, 330.34, 331.14, 331.94, 332.74, 333.53, 334.33, 335.13, 335.93, 336.73, 337.53, 338.33, 339.12, 339.92, 340.72, 341.52, 342.32, 343.11, 343.91, 344.71, 345.51, 346.3, 347.1, 347.9, 348.69, 349.49, 350.29, 351.08, 351.88, 352.67, 353.47, 354.26, 355.06, 355.85, 356.65, 357.44, 358.24, 359.03, 359.83, 360.62, 361.42, 362.21, 363, 363.8, 364.59, 365.38, 366.18, 366.97, 367.76, 368.56, 369.35, 370.14, 370.93, 371.72, 372.51, 373.31, 374.1, 374.89, 375.68, 376.47, 377.26, 378.05, 378.84, 379.63, 380.42, 381.21, 382, 382.78, 383.57, 384.36, 385.15, 385.94, 386.72, 387.51, 388.3, 389.09, 389.87, 390.66, 391.44, 392.23, 393.02, 393.8]"

# Row 8: CC_eno3 - array of coefficients (stored as text)
$ws.Range("H8").Value = "[-0.00533091, 0.00058257, 0.0089155, -0.00314105, 0.00904059, 0.01177176, 0.00505318, -0.0022058, 0.01345313, 0.00556597, -0.00232221, 0.01293107, 0.01137997, -0.0026525, -0.00013319, -0.00344627, -0.0022024, -0.00151342, -0.00839776, 0.00053789, 0.0060395, 0.00724971, 0.00904479, 0.00823721, 0.00741827, 0.0074006, 0.00695111, 0.00667655, 0.00635725, 0.00602938, 0.0056984, 0.00537362, 0.00505218, 0.00471603, 0.004395, 0.00406662, 0.00374443, 0.0034551, 0.00315349, 0.00287363, 0.00259059, 0.00232713, 0.00207377, 0.00184624, 0.00162108, 0.00142498, 0.00124222, 0.00107292, 0.00092159, 0.00078441, 0.00066109, 0.00055898, 0.00046856, 0.00038572, 0.00032082, 0.00025489, 0.00020933, 0.00016675, 0.00012926, 0.00010426, 0.00008253, 0.0000614, 0.00004665, 0.00003044, 0.00002103, 0.00000744, 0.00000078, -0.00000332, -0.00000671, -0.00001353, -0.00001482, -0.00001611, -0.00001072, -0.00001251, -0.00000883, -0.00000336, -0.00000874, 0.00000227, 0.00000567, -0.00000251, 0.00000405, -0.00000096, 0.00000323, 0.00000112, 0.00000704, 0.00000761, 0.00000756, 0.00000389, 0.00001013, 0.00001454, 0.00002335, 0.00002771, 0.00004428, 0.00003558, 0.00003293, 0.00003627, 0.0000342, 0.00003062, 0.00003455, 0.0000427, 0.00004826, 0.00005187, 0.00005684, 0.00005522, 0.00006097, 0.0000703, 0.00007116, 0.00007113, 0.00007314, 0.00008023, 0.00007779, 0.00007378, 0.0000813, 0.00007801, 0.00008188, 0.0000903, 0.00008943, 0.00010039, 0.00010451, 0.00009631, 0.00009861, 0.00010207, 0.00009975, 0.00011719, 0.00011772, 0.00012291, 0.00013072, 0.00012316, 0.00013842, 0.00012788, 0.00013165, 0.00013307, 0.00014465, 0.00014966, 0.00015677, 0.00015994, 0.00016184, 0.00015582, 0.00016465, 0.00017166, 0.0001549, 0.00017102, 0.00016414, 0.00017467, 0.0001735, 0.00017735, 0.00018267, 0.00018567, 0.0001873, 0.00019011, 0.00020054, 0.00019888, 0.00019526, 0.00019959, 0.00019569, 0.00021086, 0.00020261, 0.00021149, 0.00021463, 0.00021403, 0.00022225, 0.00022225, 0.00022285, 0.00022626, 0.00023876, 0.00023162, 0.00023825, 0.00023202, 0.00024257, 0.00024408, 0.00024866, 0.00025064, 0.0002489, 0.0002518, 0.00027095, 0.00027383, 0.00026756, 0.00028019, 0.00026779, 0.00027804, 0.00027763, 0.00028107, 0.00028621, 0.00027606, 0.00028923, 0.00028403, 0.00028393, 0.00029102, 0.00029822, 0.00029578, 0.00029818, 0.0003092, 0.00030734, 0.00031047, 0.00032424, 0.00031928, 0.00031999, 0.00032696, 0.00032952, 0.00031152, 0.00032103, 0.00033113, 0.00033504, 0.00033148, 0.00033907, 0.00034968, 0.00034129, 0.00034333, 0.00034598, 0.00035525, 0.00035845, 0.00035981, 0.00036662, 0.00036611, 0.00038937, 0.00038152, 0.00039393, 0.00037037, 0.00037274, 0.00037765, 0.00038838, 0.00037269, 0.00038598, 0.0004036, 0.00041153, 0.00040715, 0.00042598, 0.00041857, 0.00041965, 0.00040934, 0.00041704, 0.0004281, 0.00042906, 0.00043078, 0.00042345, 0.00043417, 0.00041648, 0.00044515, 0.00043906, 0.00043705, 0.0004448, 0.00045549, 0.00043654, 0.00045781, 0.0004449, 0.00045409, 0.00048519, 0.00046696, 0.00047051, 0.00045589, 0.00046137, 0.00047814, 0.00047453, 0.0004801, 0.00045135, 0.00047998]"

# Row 9: CC_eswa - array of coefficients (stored as text)
$ws.Range("H9").Value = "[0.00619525, -0.00132526, 0.00084339, 0.00091014, -0.00088382, 0.00434245, 0.00873569, 0.00531021, 0.0137682, 0.01238712, 0.03094334, 0.04514063, 0.05480445, 0.06763339, 0.07699255, 0.0805132, 0.09032271, 0.09143394, 0.09918138, 0.08842984, 0.07986757, 0.0708021, 0.06209281, 0.05344975, 0.04539421, 0.03799354, 0.03163095, 0.02606925, 0.02133065, 0.01734092, 0.01402406, 0.01127976, 0.0090108, 0.00718585, 0.00570159, 0.00451514, 0.00356472, 0.00281227, 0.00221345, 0.00173197, 0.00136032, 0.0010719, 0.00085594, 0.00066455, 0.00053052, 0.00042122, 0.00033881, 0.00026315, 0.00023223, 0.00018588, 0.00016216, 0.00013129, 0.00010911, 0.00010268, 0.00007923, 0.00007073, 0.00005857, 0.00005558, 0.00004677, 0.00004417, 0.00002835, 0.00002746, 0.00001805, 0.00001592, 0.00000546, 0.00000399, 0.00000413, -0.00000004, 0.00000493, 0.00000749, -0.00000063, -0.0000016, -0.00000427, -0.00000972, -0.00001194, -0.00001295, 0.00000093, -0.0000126, -0.00001038, -0.0000118, -0.00000591, -0.00000759, -0.0000099, 0.00000241, 0.00001084, 0.00001132, 0.00001306, 0.00001885, 0.00001539, 0.00002456, 0.00001455, 0.00001942, 0.00001054, 0.00001866, 0.00002203, 0.00002297, 0.000035, 0.00003712, 0.00003105, 0.00003569, 0.00004276, 0.00005292, 0.00004396, 0.00004979, 0.00005409, 0.00005085, 0.00005717, 0.00005154, 0.00005848, 0.00006412, 0.00006, 0.00006987, 0.00006794, 0.00006121, 0.00007415, 0.0000693, 0.00007832, 0.00008025, 0.00008023, 0.00008622, 0.00008201, 0.00008401, 0.00009803, 0.0000826, 0.0000896, 0.00007205, 0.00009118, 0.0000937, 0.00009561, 0.00010279, 0.00009604, 0.00011781, 0.00010901, 0.00010615, 0.00010516, 0.00011212, 0.00010695, 0.00012016, 0.00011828, 0.00012683, 0.00013331, 0.00012705, 0.00013528, 0.00013049, 0.0001374, 0.00014344, 0.00014779, 0.00014635, 0.00015073, 0.0001483, 0.00014928, 0.00015925, 0.00017176, 0.00016701, 0.00017956, 0.00017731, 0.00017254, 0.00017577, 0.00017092, 0.0001776, 0.00017351, 0.00018404, 0.00019135, 0.00019609, 0.00019996, 0.00020709, 0.0002018, 0.00022587, 0.00021064, 0.00021278, 0.00021918, 0.00023631, 0.00023022, 0.00022719, 0.0002418, 0.00022582, 0.00024021, 0.00023473, 0.00024183, 0.00025958, 0.00025052, 0.00024722, 0.00026701, 0.00026712, 0.00027277, 0.00027781, 0.00028209, 0.00028941, 0.00027998, 0.00030498, 0.00029836, 0.00029083, 0.00030658, 0.00031418, 0.00030212, 0.00031807, 0.00031969, 0.00034033, 0.00032975, 0.0003648, 0.00034142, 0.00034686, 0.00034859, 0.000354, 0.00034526, 0.00032918, 0.00036067, 0.00035943, 0.00037007, 0.00035261, 0.00038674, 0.00037957, 0.00037093, 0.00037763, 0.00037622, 0.00039027, 0.00039695, 0.00042034, 0.00041343, 0.00041097, 0.00041934, 0.00041431, 0.00042816, 0.00041475, 0.00041951, 0.0004192, 0.00039958, 0.00041802, 0.00043268, 0.00044414, 0.00044557, 0.00045845, 0.00045221, 0.00048219, 0.00046398, 0.00046445, 0.00047098, 0.00045474, 0.0004683, 0.00048103, 0.00047327, 0.00045607, 0.00050671, 0.00046632, 0.00051272, 0.00052347, 0.00048511, 0.00052727, 0.00052736, 0.00052374, 0.00053235, 0.00053614, 0.00053242, 0.00051571, 0.00055037, 0.00055533]"

# Row 10: CC_di - array of coefficients (stored as text)
$ws.Range("H10").Value = "[21.58333333, 26.16666667, 25.5, 27.375, 23.66666667, 29.29166667, 27.5, 23.70833333, 22.20833333, 51.04166667, 186.54166667, 641.125, 1782.70833333, 4117.375, 7972.29166667, 13085.625, 18466.20833333, 22998.66666667, 26124.83333333, 27859.16666667, 28524.91666667, 28572.875, 28361.41666667, 28139.16666667, 28120.5, 28359.20833333, 28865.41666667, 29741.875, 30941.79166667, 32520.20833333, 34432.66666667, 36718.91666667, 39286.54166667, 41951.375, 44608, 47108.25, 49282.79166667, 50866.58333333, 51668.125, 51685.95833333, 50807.125, 49291.41666667, 47295.16666667, 45028.75, 42750.45833333, 40640.08333333, 38781.95833333, 37284.16666667, 36157.16666667, 35386.5, 34955.16666667, 34885.04166667, 35084.25, 35568.58333333, 36319.66666667, 37318.5, 38472.16666667, 39853.20833333, 41376.25, 43089.875, 44871.625, 46670.91666667, 48400.70833333, 50037.45833333, 51429.45833333, 52467.91666667, 53088.125, 53132.45833333, 52671.54166667, 51628.70833333, 50093.08333333, 48228.58333333, 46068.125, 43726.83333333, 41414.04166667, 39114.125, 36957.29166667, 34997.29166667, 33235.08333333, 31652.54166667, 30282.66666667, 29095.625, 28092.91666667, 27293.91666667, 26616.54166667, 26101.54166667, 25756.66666667, 25546.25, 25472.54166667, 25564.875, 25774.04166667, 26120.25, 26592.95833333, 27176.91666667, 27878.58333333, 28667.75, 29480.83333333, 30323.58333333, 31143.625, 31912.5, 32640.125, 33177.16666667, 33552.375, 33745.04166667, 33720.25, 33446.33333333, 32986.75, 32328.70833333, 31521.25, 30597.54166667, 29577.125, 28522.33333333, 27489.5, 26450.875, 25529.95833333, 24673.54166667, 23891.45833333, 23189.5, 22601.04166667, 22086.33333333, 21672.29166667, 21366.45833333, 21140.66666667, 21014.08333333, 20975.83333333, 21002.375, 21132.66666667, 21328.66666667, 21583.29166667, 21933.875, 22316.41666667, 22775.04166667, 23277.29166667, 23834.54166667, 24436.83333333, 25093.625, 25783.29166667, 26483.33333333, 27213.375, 27947, 28695.70833333, 29452.33333333, 30215.75, 30934.08333333, 31619, 32231.54166667, 32780.125, 33289.58333333, 33658.08333333, 33929.75, 34095.08333333, 34138.5, 34073.75, 33928.16666667, 33623.75, 33262.20833333, 32797.375, 32256.45833333, 31701.04166667, 31081.58333333, 30418.875, 29804.20833333, 29155, 28478.75, 27782.04166667, 27088.25, 26334.41666667, 25621.41666667, 24909.58333333, 24227.125, 23571.125, 22953.125, 22370.16666667, 21846.45833333, 21356.41666667, 20923.79166667, 20529.875, 20145.25, 19796.16666667, 19485.58333333, 19178.75, 18912.625, 18666.125, 18421.83333333, 18183.75, 17982.91666667, 17774.625, 17610.75, 17459.79166667, 17306.33333333, 17169, 17056.29166667, 16962.91666667, 16900.625, 16858.875, 16828.70833333, 16822.70833333, 16831.70833333, 16838, 16858.25, 16896.29166667, 16949.70833333, 17018.58333333, 17076.54166667, 17137.54166667, 17205.08333333, 17250.95833333, 17311.375, 17359.54166667, 17390.54166667, 17410.04166667, 17418.125, 17402.75, 17348.08333333, 17280.08333333, 17172.625, 17050.5, 16914.83333333, 16774.5, 16626.5, 16470.91666667, 16283.04166667, 16065.45833333, 15797.70833333, 15469.54166667, 15132.91666667, 14785.04166667, 14460.54166667, 14177, 13903.58333333, 13650.25, 13425.66666667, 13193.58333333, 13004.45833333, 12835.45833333, 12674.75, 12521.75, 12354.5, 12137.54166667, 11890.95833333, 11656.54166667, 11427.54166667, 11243.58333333, 11069.16666667, 10913, 10785.04166667, 10666, 10568.83333333, 10487.75, 10429.41666667, 10382.83333333, 10345.95833333, 10314.75, 10177.16666667, 9821, 9019.54166667]"

# Row 11: CC_lower_wavelength_limit_for_spectra_fit - numeric, new highlighted style
$ws.Range("I11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = 217
$ws.Range("H11").Interior.Color = 5296274
$ws.Range("H11").NumberFormat = "0.00"

# Row 12: CC_upper_wavelength_limit_for_spectra_fit - numeric, new highlighted style
$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = 240

# Selection state: move selection to F17 on Asset_Cal_Info
[void]$ws.Range("F17").Select()

# Switch active sheet back to Moorings
$ws1 = $wb.Worksheets.Item("Moorings")
[void]$ws1.Select()
